$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy header formatting (bold / border / centered) from an existing
#     header cell onto the new header cells L1:O1 before filling values ---
$ws.Range("K1").Copy()
$ws.Range("L1:O1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Header row (row 1) ---
$ws.Range("G1").Value = "점수(룰)"
$ws.Range("H1").Value = "3일상승확률(%)"
$ws.Range("I1").Value = "5일상승확률(%)"
$ws.Range("J1").Value = "10일상승확률(%)"
$ws.Range("K1").Value = "최종점수"
$ws.Range("L1").Value = "예측방식"
$ws.Range("M1").Value = "판단"
$ws.Range("N1").Value = "MACRO_SCORE"
$ws.Range("O1").Value = "MACRO_SIGNAL"

# --- Row 2 (now Oklo Inc. / OKLO) ---
$ws.Range("A2").Value = "'2025-11-29"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "Oklo Inc."
$ws.Range("C2").Value = "OKLO"
$ws.Range("D2").Value = 91.88
$ws.Range("E2").ClearContents()
$ws.Range("F2").Value = 4.41
$ws.Range("G2").Value = 20
$ws.Range("H2").Value = 70
$ws.Range("I2").Value = 76
$ws.Range("J2").Value = 83
$ws.Range("K2").Value = 62
$ws.Range("L2").Value = "Pattern"
$ws.Range("M2").Value = "📈 매수 관찰 구간입니다."
$ws.Range("N2").Value = 85.36763896678245
$ws.Range("O2").Value = "🟢 완화적 (상승 우위)"

# --- Row 3 (now NuScale Power Corporation / SMR) ---
$ws.Range("A3").Value = "'2025-11-29"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "NuScale Power Corporation"
$ws.Range("C3").Value = "SMR"
$ws.Range("D3").Value = 20
$ws.Range("E3").ClearContents()
$ws.Range("F3").Value = 6.95
$ws.Range("G3").Value = 20
$ws.Range("H3").Value = 63
$ws.Range("I3").Value = 60
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 55.6
$ws.Range("L3").Value = "Pattern"
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 85.36763896678245
$ws.Range("O3").Value = "🟢 완화적 (상승 우위)"
